$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("G2").Value = 92
$ws.Range("I2").Value = 120
$ws.Range("F3").Value = 147
$ws.Range("E3").Value = 152
$ws.Range("G3").Value = 149
$ws.Range("H6").Value = 481
$ws.Range("F6").Value = 588
$ws.Range("D6").Value = 444
$ws.Range("E6").Value = 514
$ws.Range("G6").Value = 449
$ws.Range("F7").Value = 843
$ws.Range("G7").Value = 696
$ws.Range("D7").Value = 687
$ws.Range("I7").Value = 864
$ws.Range("H7").Value = 770
$ws.Range("E7").Value = 753

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("F5").Value = 15
$ws.Range("F6").Value = 19

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("F3").Value = 1
$ws.Range("F5").Value = 25
$ws.Range("F6").Value = 28

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("I2").Value = 1
$ws.Range("I5").Value = 12

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("F6").Value = 42
$ws.Range("F7").Value = 63

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("F5").Value = 19
$ws.Range("F28").Value = 63
$ws.Range("G30").Value = 5
$ws.Range("H47").Value = 27
$ws.Range("F50").Value = 28
$ws.Range("F53").Value = 91
$ws.Range("G53").Value = 88
$ws.Range("E53").Value = 91
$ws.Range("E65").Value = 21
$ws.Range("G76").Value = 20
$ws.Range("F78").Value = 13
$ws.Range("E78").Value = 11
$ws.Range("D81").Value = 4
$ws.Range("I88").Value = 12
$ws.Range("F91").Value = 12
$ws.Range("F98").Value = 843
$ws.Range("E98").Value = 753
$ws.Range("H98").Value = 770
$ws.Range("G98").Value = 696
$ws.Range("D98").Value = 687
$ws.Range("I98").Value = 864

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("E3").Value = 3
$ws.Range("E4").Value = 7
$ws.Range("F4").Value = 9
$ws.Range("F5").Value = 13
$ws.Range("E5").Value = 11

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("G3").Value = 26
$ws.Range("F6").Value = 68
$ws.Range("E6").Value = 67
$ws.Range("F7").Value = 91
$ws.Range("G7").Value = 88
$ws.Range("E7").Value = 91

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("F6").Value = 10
$ws.Range("F7").Value = 12

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("D5").Value = 3
$ws.Range("D6").Value = 4

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("G2").Value = 5
$ws.Range("G7").Value = 20

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("E5").Value = 18
$ws.Range("E6").Value = 21

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("G6").Value = 4
$ws.Range("G7").Value = 5

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("H5").Value = 17
$ws.Range("H6").Value = 27
